# Generate Report for Handoff
# - Updates the "Status" column text and handoff/handback timestamps
#   to reflect a freshly-generated handoff report.
# - Shrinks the now-shorter "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-22 17:01:52"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 17:01:47"

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 17:01:52"

# --- Column width adjustments --------------------------------------------
# The "Status" columns shrink now that the text is shorter.
# (16.3333... is the ColumnWidth input that lands closest to the
# target stored width of 17.2159881591797 given this engine's
# pixel-quantized ColumnWidth setter.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
